$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    $p = $d.Paragraphs.Item($Index)
    $r = $p.Range
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $found) {
        Write-Output ("WARNING: not found in paragraph " + $Index + ": " + $OldText)
    }
}

$BR = [char]11

# --- Paragraph 6: "Objetivos" body (PT) ---
$old06 = "Estudo de Óptica Física."
$new06 = "Descrição ondulatória e quântica da luz. Propriedades da luz. Interação da luz com a matéria. Aplicações."
Replace-InParagraph 6 $old06 $new06

# --- Paragraph 7: "Objetivos" body (EN, italic) ---
$old07 = "Study of Physical Optics."
$new07 = "Presentation of the wave and quantum description of light, study of the properties of light, the interaction of light with matter and applications of physical optics."
Replace-InParagraph 7 $old07 $new07

# --- Paragraph 9: "Docente(s) Responsável(eis)" bullet list (two runs) ---
$old09a = "519033 - Carlos Yujiro Shigue"
$new09a = "Estudo de Óptica Física."
Replace-InParagraph 9 $old09a $new09a

$old09b = "1643715 - Paulo Atsushi Suzuki"
$new09b = "O que é luz? Reflexão. Refração. Difração. Polarização. Formação de imagens: Transformada de Fourier. Ondas eletromagnéticas. Equações de Maxwell. Propagação da luz em diferentes meios: vácuo, dielétrico, condutor. Transporte de energia. Condições de contorno entre diferentes meios: vácuo, dielétrico, condutor. Propagação da luz entre diferentes meios: incidência normal e oblíqua na interface entre meios. Coeficientes de Fresnel. Aplicações da Óptica: holografia, laser, fibras ópticas, materiais eletrocrômicos, metamateriais."
Replace-InParagraph 9 $old09b $new09b

# --- Paragraph 11: "Programa resumido" body (PT) ---
$old11 = "Descrição ondulatória e quântica da luz. Propriedades da luz. Interação da luz com a matéria. Aplicações."
$new11 = "Aulas expositivas, resolução de exercícios e seminários."
Replace-InParagraph 11 $old11 $new11

# --- Paragraph 12: "Programa resumido" body (EN, italic) ---
$old12 = "Presentation of the wave and quantum description of light, study of the properties of light, the interaction of light with matter and applications of physical optics."
$new12 = "Study of Physical Optics."
Replace-InParagraph 12 $old12 $new12

# --- Paragraph 14: "Programa" body (PT) ---
$old14 = "O que é luz? Reflexão. Refração. Difração. Polarização. Formação de imagens: Transformada de Fourier. Ondas eletromagnéticas. Equações de Maxwell. Propagação da luz em diferentes meios: vácuo, dielétrico, condutor. Transporte de energia. Condições de contorno entre diferentes meios: vácuo, dielétrico, condutor. Propagação da luz entre diferentes meios: incidência normal e oblíqua na interface entre meios. Coeficientes de Fresnel. Aplicações da Óptica: holografia, laser, fibras ópticas, materiais eletrocrômicos, metamateriais."
$new14 = "Média ponderada de duas provas escritas: P1, P2 e TR. Conceito Final = (P1 + 2P2)/3"
Replace-InParagraph 14 $old14 $new14

# --- Paragraph 17: "Avaliação" bullet (Método / Critério / Norma de recuperação) ---
# IMPORTANT: do the "Norma de recuperação:" content swap BEFORE the "Método:" content
# swap, because the new "Método:" text equals the old "Norma de recuperação:" text --
# replacing Método first would create a duplicate that the later Find could mismatch.
$old17a = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$new17a = "519033 - Carlos Yujiro Shigue"
Replace-InParagraph 17 $old17a $new17a

$old17b = "Média ponderada de duas provas escritas: P1, P2 e TR. Conceito Final = (P1 + 2P2)/3"
$new17b = "HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974." + $BR + "ZILLIO, S. C. Óptica Moderna - Fundamentos e Aplicações, 2005." + $BR + "J. R. Reitz, F. J. Milford, R. W. Christy, Fundamentos da Teoria Eletromagnética. Editora Campus. 1982."
Replace-InParagraph 17 $old17b $new17b

$old17c = "Aulas expositivas, resolução de exercícios e seminários."
$new17c = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
Replace-InParagraph 17 $old17c $new17c

# --- Paragraph 19: "Bibliografia" body (three lines -> single line) ---
$old19 = "HECHT, E.; ZAJAC, A. Optics; Reading, Addison-Wesley, 1974." + $BR + "ZILLIO, S. C. Óptica Moderna - Fundamentos e Aplicações, 2005." + $BR + "J. R. Reitz, F. J. Milford, R. W. Christy, Fundamentos da Teoria Eletromagnética. Editora Campus. 1982."
$new19 = "1643715 - Paulo Atsushi Suzuki"
Replace-InParagraph 19 $old19 $new19

Write-Output "Done."
